# Weekly update: a new week of price data (2022-08-25, serial 44798) is
# prepended to the Betarraga price history. This pushes the existing
# rows 274-325 down by two rows (to 276-327); the two freed rows
# (274-275) are filled with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 274, shifting rows
# 274:325 down to 276:327 (same effect as the data cascading down
# one "week" to make room for the newest week at the top).
$ws.Rows("274:275").Insert()

# New first data row of the newly-inserted week (Calidad = Primera)
$ws.Range("A274").Value = 8
$ws.Range("B274").Value = "Terminal La Palmera de La Serena"
$ws.Range("C274").Value = "Coquimbo"
$ws.Range("D274").Value = 44798
$ws.Range("E274").Value = 4
$ws.Range("F274").Value = 100114014
$ws.Range("G274").Value = "Betarraga"
$ws.Range("H274").Value = "Sin especificar"
$ws.Range("I274").Value = "Primera"
$ws.Range("J274").Value = 2000
$ws.Range("K274").Value = 600
$ws.Range("L274").Value = 700
$ws.Range("M274").Value = 650
$ws.Range("N274").Value = "$/paquete 3 unidades"
$ws.Range("O274").Value = "Provincia del Elquí"
$ws.Range("P274").Value = 217
$ws.Range("Q274").Value = 3
$ws.Range("R274").Value = "Hortaliza"

# New second data row of the newly-inserted week (Calidad = Segunda)
$ws.Range("A275").Value = 8
$ws.Range("B275").Value = "Terminal La Palmera de La Serena"
$ws.Range("C275").Value = "Coquimbo"
$ws.Range("D275").Value = 44798
$ws.Range("E275").Value = 4
$ws.Range("F275").Value = 100114014
$ws.Range("G275").Value = "Betarraga"
$ws.Range("H275").Value = "Sin especificar"
$ws.Range("I275").Value = "Segunda"
$ws.Range("J275").Value = 1520
$ws.Range("K275").Value = 500
$ws.Range("L275").Value = 550
$ws.Range("M275").Value = 525
$ws.Range("N275").Value = "$/paquete 3 unidades"
$ws.Range("O275").Value = "Provincia del Elquí"
$ws.Range("P275").Value = 175
$ws.Range("Q275").Value = 3
$ws.Range("R275").Value = "Hortaliza"
